$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# --- Add PhoneNumber column (J) to the Gmail_Signup_Testdata sheet ---
$ws1.Range("J1").Value = "PhoneNumber"
$ws1.Range("J2").Value = "6504603326"
$ws1.Columns.Item(10).ColumnWidth = 13.26953125

# clear the stray selection left on sheet1 and select A1:D3 instead
$ws1.Range("A1:D3").Select()

# --- Add the new Gmail_Signin_Data sheet, placed after the existing sheet ---
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Gmail_Signin_Data"

$ws2.Range("A1").Value = "Scenario"
$ws2.Range("B1").Value = "TestCase"
$ws2.Range("C1").Value = "UserEmail"
$ws2.Range("D1").Value = "Password"
$ws2.Range("E1").Value = "FirstName"
$ws2.Range("F1").Value = "LastName"

$ws2.Range("A2").Value = "User1"
$ws2.Range("B2").Value = "1"
$ws2.Range("C2").Value = "vinayaknaiktest1@gmail.com"
$ws2.Range("D2").Value = "January@123"
$ws2.Range("E2").Value = "Vinayak"
$ws2.Range("F2").Value = "Naik"

$ws2.Range("A3").Value = "User2"
$ws2.Range("B3").Value = "2"

$ws2.Hyperlinks.Add($ws2.Range("C2"), "mailto:vinayaknaiktest1@gmail.com")
$ws2.Hyperlinks.Add($ws2.Range("D2"), "mailto:January@123")

$ws2.Columns.Item(3).ColumnWidth = 25.36328125
$ws2.Columns.Item(4).ColumnWidth = 13
$ws2.Columns.Item(5).ColumnWidth = 12.6328125
$ws2.Columns.Item(6).ColumnWidth = 13

$ws2.Range("E2:F3").Select()
